$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# D5: 34 -> 33
$wsForecast.Range("D5").Value = 33

# D14: 25 -> 26
$wsForecast.Range("D14").Value = 26

# D15: 26 -> 25
$wsForecast.Range("D15").Value = 25

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")

# Force these cells to remain plain text, matching the original inline-string
# storage, instead of letting Excel auto-convert numeric/date-looking text.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B15").NumberFormat = "@"

# B9: "502" -> "501"
$wsSummary.Range("B9").Value = "501"

# B10: "266" -> "265"
$wsSummary.Range("B10").Value = "265"

# B11: "116" -> "115"
$wsSummary.Range("B11").Value = "115"

# B15: "2025-04-20" -> "2025-04-27"
$wsSummary.Range("B15").Value = "2025-04-27"
